$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()

$ws.Range("H40").Value = 4180.8887
$ws.Range("J40").Value = 4592.4546
$ws.Range("L40").Value = 4592.4546
$ws.Range("N40").Value = -4942.4546

$ws.Range("H132").Value = 11216.431
$ws.Range("I132").Value = 1893.5652
$ws.Range("J132").Value = 46954.082
$ws.Range("K132").Value = 5680.6956
$ws.Range("L132").Value = 140862.246
$ws.Range("M132").Value = -3150.6956
$ws.Range("N132").Value = -145922.246

$ws.Range("H135").Value = 27780518
$ws.Range("J135").Value = 3012
$ws.Range("L135").Value = 27108
$ws.Range("N135").Value = -32178

$ws.Range("H137").Value = 2444.0588
$ws.Range("I137").Value = 2169.9666
$ws.Range("J137").Value = 4499.75
$ws.Range("K137").Value = 6509.899800000001
$ws.Range("L137").Value = 13499.25
$ws.Range("M137").Value = -3959.899800000001
$ws.Range("N137").Value = -18599.25

$ws.Range("H139").Value = 90000
$ws.Range("J139").Value = 90000
$ws.Range("L139").Value = 90000
$ws.Range("N139").Value = -100280

$ws.Range("H141").Value = 5039.16
$ws.Range("I141").Value = 3958.524
$ws.Range("K141").Value = 11875.572
$ws.Range("M141").Value = -6695.572

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3534.7856
$ws.Range("I2").Value = 3347.4
$ws.Range("K2").Value = 3347.4
$ws.Range("M2").Value = -3234.4

$ws.Range("H32").Value = 2908.6206
$ws.Range("I32").Value = 2898.077
$ws.Range("J32").Value = 3000
$ws.Range("K32").Value = 2898.077
$ws.Range("L32").Value = 3000
$ws.Range("M32").Value = -2611.077
$ws.Range("N32").Value = -3574

$ws.Range("H61").Value = 2669.7
$ws.Range("I61").Value = 2669.7
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2669.7
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2457.7
$ws.Range("N61").ClearContents()

$ws.Range("H74").Value = 1995.125
$ws.Range("I74").Value = 2093.1
$ws.Range("J74").Value = 1505.25
$ws.Range("K74").Value = 2093.1
$ws.Range("L74").Value = 1505.25
$ws.Range("M74").Value = -1219.1
$ws.Range("N74").Value = -3253.25

$ws.Range("H77").Value = 1995.125
$ws.Range("I77").Value = 2093.1
$ws.Range("J77").Value = 1505.25
$ws.Range("K77").Value = 10465.5
$ws.Range("L77").Value = 7526.25
$ws.Range("M77").Value = -6097.5
$ws.Range("N77").Value = -16262.25

$ws.Range("H116").Value = 3534.7856
$ws.Range("I116").Value = 3347.4
$ws.Range("K116").Value = 3347.4
$ws.Range("M116").Value = -1053.4

$ws.Range("H132").Value = 1744.625
$ws.Range("I132").Value = 1744.625
$ws.Range("K132").Value = 5233.875
$ws.Range("M132").Value = -2703.875

$ws.Range("H136").Value = 2669.7
$ws.Range("I136").Value = 2669.7
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 8009.099999999999
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -5459.099999999999
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3534.7856
$ws.Range("I3").Value = 3347.4
$ws.Range("K3").Value = 3347.4
$ws.Range("M3").Value = -3233.4

$ws.Range("H86").Value = 6397.3335
$ws.Range("I86").Value = 3676.8
$ws.Range("J86").Value = 20000
$ws.Range("K86").Value = 3676.8
$ws.Range("L86").Value = 20000
$ws.Range("M86").Value = -2553.8
$ws.Range("N86").Value = -22246

$ws.Range("H89").Value = 6397.3335
$ws.Range("I89").Value = 3676.8
$ws.Range("J89").Value = 20000
$ws.Range("K89").Value = 18384
$ws.Range("L89").Value = 100000
$ws.Range("M89").Value = -12768
$ws.Range("N89").Value = -111232

$ws.Range("H107").Value = 7574.2666
$ws.Range("I107").Value = 6758.4165
$ws.Range("J107").Value = 10837.667
$ws.Range("K107").Value = 6758.4165
$ws.Range("L107").Value = 10837.667
$ws.Range("M107").Value = -4838.4165
$ws.Range("N107").Value = -14677.667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2517.0476
$ws.Range("I16").Value = 2387.3157
$ws.Range("J16").Value = 3749.5
$ws.Range("K16").Value = 2387.3157
$ws.Range("L16").Value = 3749.5
$ws.Range("M16").Value = -2100.3157
$ws.Range("N16").Value = -4323.5

$ws.Range("H31").Value = 1297.2858
$ws.Range("I31").Value = 1232.6364
$ws.Range("J31").Value = 1534.3334
$ws.Range("K31").Value = 1232.6364
$ws.Range("L31").Value = 1534.3334
$ws.Range("M31").Value = -937.6364000000001
$ws.Range("N31").Value = -2124.3334

$ws.Range("H34").Value = 1297.2858
$ws.Range("I34").Value = 1232.6364
$ws.Range("J34").Value = 1534.3334
$ws.Range("K34").Value = 1232.6364
$ws.Range("L34").Value = 1534.3334
$ws.Range("M34").Value = -1030.6364
$ws.Range("N34").Value = -1938.3334

$ws.Range("H58").Value = 2857.1177
$ws.Range("I58").Value = 1709
$ws.Range("J58").Value = 3877.6667
$ws.Range("K58").Value = 1709
$ws.Range("L58").Value = 3877.6667
$ws.Range("M58").Value = -1506
$ws.Range("N58").Value = -4283.6667

$ws.Range("H99").Value = 3547.923
$ws.Range("J99").Value = 3284.1875
$ws.Range("L99").Value = 3284.1875
$ws.Range("N99").Value = -6280.1875

$ws.Range("H107").Value = 6581.222
$ws.Range("I107").Value = 1134
$ws.Range("J107").Value = 12028.444
$ws.Range("K107").Value = 1134
$ws.Range("L107").Value = 12028.444
$ws.Range("M107").Value = 786
$ws.Range("N107").Value = -15868.444

$ws.Range("H113").Value = 2517.0476
$ws.Range("I113").Value = 2387.3157
$ws.Range("J113").Value = 3749.5
$ws.Range("K113").Value = 2387.3157
$ws.Range("L113").Value = 3749.5
$ws.Range("M113").Value = -217.3157000000001
$ws.Range("N113").Value = -8089.5

$ws.Range("H126").Value = 3547.923
$ws.Range("J126").Value = 3284.1875
$ws.Range("L126").Value = 9852.5625
$ws.Range("N126").Value = -14792.5625

$ws.Range("H132").Value = 2776.8518
$ws.Range("I132").Value = 2000.3334
$ws.Range("J132").Value = 4329.8887
$ws.Range("K132").Value = 6001.0002
$ws.Range("L132").Value = 12989.6661
$ws.Range("M132").Value = -3471.0002
$ws.Range("N132").Value = -18049.6661

$ws.Range("H134").Value = 5102.5
$ws.Range("I134").Value = 5106.3477
$ws.Range("J134").Value = 5014
$ws.Range("K134").Value = 15319.0431
$ws.Range("L134").Value = 15042
$ws.Range("M134").Value = -12784.0431
$ws.Range("N134").Value = -20112

$ws.Range("H136").Value = 2857.1177
$ws.Range("I136").Value = 1709
$ws.Range("J136").Value = 3877.6667
$ws.Range("K136").Value = 5127
$ws.Range("L136").Value = 11633.0001
$ws.Range("M136").Value = -2577
$ws.Range("N136").Value = -16733.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 71.25
$ws.Range("I2").Value = 69.5
$ws.Range("J2").Value = 73
$ws.Range("K2").Value = 417
$ws.Range("L2").Value = 438
$ws.Range("M2").Value = -304
$ws.Range("N2").Value = -664

$ws.Range("H38").Value = 1521.2142
$ws.Range("I38").Value = 91.166664
$ws.Range("J38").Value = 2593.75
$ws.Range("K38").Value = 273.499992
$ws.Range("L38").Value = 7781.25
$ws.Range("M38").Value = 73.50000799999998
$ws.Range("N38").Value = -8475.25

$ws.Range("H95").Value = 6633.3335
$ws.Range("J95").Value = 6633.3335
$ws.Range("L95").Value = 19900.0005
$ws.Range("N95").Value = -24018.0005

$ws.Range("H132").Value = 2927.0908
$ws.Range("I132").Value = 3179.6
$ws.Range("K132").Value = 28616.4
$ws.Range("M132").Value = -26086.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 10604.35
$ws.Range("I132").Value = 9828.647000000001
$ws.Range("J132").Value = 15000
$ws.Range("K132").Value = 29485.941
$ws.Range("L132").Value = 45000
$ws.Range("M132").Value = -26955.941
$ws.Range("N132").Value = -50060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3778.4
$ws.Range("I7").Value = 2868.3845
$ws.Range("K7").Value = 2868.3845
$ws.Range("M7").Value = -2756.3845

$ws.Range("H16").Value = 880
$ws.Range("I16").Value = 758.75
$ws.Range("J16").Value = 1041.6666
$ws.Range("K16").Value = 758.75
$ws.Range("L16").Value = 1041.6666
$ws.Range("M16").Value = -588.75
$ws.Range("N16").Value = -1381.6666

$ws.Range("H40").Value = 2696.4666
$ws.Range("I40").Value = 2107.25
$ws.Range("K40").Value = 2107.25
$ws.Range("M40").Value = -1971.25

$ws.Range("H61").Value = 2762.125
$ws.Range("I61").Value = 2728.2856
$ws.Range("K61").Value = 2728.2856
$ws.Range("M61").Value = -2526.2856

$ws.Range("H113").Value = 2762.125
$ws.Range("I113").Value = 2728.2856
$ws.Range("K113").Value = 2728.2856
$ws.Range("M113").Value = -558.2856000000002

$ws.Range("H126").Value = 3778.4
$ws.Range("I126").Value = 2868.3845
$ws.Range("K126").Value = 8605.1535
$ws.Range("M126").Value = -6135.1535

$ws.Range("H132").Value = 1889.6666
$ws.Range("I132").Value = 1922.04
$ws.Range("J132").Value = 1485
$ws.Range("K132").Value = 5766.12
$ws.Range("L132").Value = 4455
$ws.Range("M132").Value = -3236.12
$ws.Range("N132").Value = -9515

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2199.75
$ws.Range("I132").Value = 1433
$ws.Range("K132").Value = 4299
$ws.Range("M132").Value = -1769
